{"js": "// Revert \"Wireframes version 2.\" -> restore \"Version 1.\"\n// Before: \"Versi\" | \"on\" | \" 2\" | \".\" (runs split around a spell-check run\n//         and the _GoBack bookmark)  => visible text \"Version 2.\"\n// After:  \"Version\" | \" 1.\" (bookmark preserved, trailing \".\" run removed)\n//         => visible text \"Version 1.\"\n\n// Step 1: merge the \"Versi\"/\"on\" runs back into a single \"Version\" run by\n// replacing that exact span with itself (search spans both runs; the\n// replace rewrites them as one run).\nlet versionResults = context.document.body.search(\"Version\", { matchCase: true });\nversionResults.load(\"text\");\nawait context.sync();\nversionResults.items[0].insertText(\"Version\", \"Replace\");\nawait context.sync();\n\n// Step 2: change the \"2\" to \"1.\" in the \" 2\" run (this run sits right\n// before the _GoBack bookmark, which stays untouched).\nlet numberResults = context.document.body.search(\" 2\", { matchCase: true });\nnumberResults.load(\"text\");\nawait context.sync();\nnumberResults.items[0].insertText(\" 1.\", \"Replace\");\nawait context.sync();\n\n// Step 3: drop the now-redundant trailing \".\" run (the one that used to\n// follow the bookmark) -- find every \".\" and remove the last occurrence.\nlet dotResults = context.document.body.search(\".\", { matchCase: true, matchWildcards: false });\ndotResults.load(\"items\");\nawait context.sync();\nif (dotResults.items.length > 1) {\n  dotResults.items[dotResults.items.length - 1].insertText(\"\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Revert \"Wireframes version 2.\" -> restore \"Version 1.\"\n# Before: \"Versi\" | \"on\" | \" 2\" | \".\" (runs split around a spell-check run\n#         and the _GoBack bookmark)  => visible text \"Version 2.\"\n# After:  \"Version\" | \" 1.\" (bookmark preserved, trailing \".\" run removed)\n#         => visible text \"Version 1.\"\n\n$d = $word.ActiveDocument\n\n# Step 1: merge the \"Versi\"/\"on\" runs back into a single \"Version\" run by\n# doing a no-op Find/Replace over that exact span (Word rewrites the match\n# as a single run).\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\"Version\", $true, $false, $false, $false, $false, $true, 1, $false, \"Version\", 2) | Out-Null\n\n# Step 2: change the \"2\" to \"1.\" in the \" 2\" run (this run sits right\n# before the _GoBack bookmark, which stays untouched by this search).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\" 2\", $false, $false, $false, $false, $false, $true, 1, $false, \" 1.\", 2) | Out-Null\n\n# Step 3: drop the now-redundant trailing \".\" run (the one that used to\n# follow the bookmark) by deleting the final character of the document.\n$tail = $d.Range($d.Content.End - 2, $d.Content.End - 1)\nif ($tail.Text -eq \".\") {\n    $tail.Delete()\n}\n"}
